$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert the new "2022-Q3" sheet right before "2022-Q2".
# ---------------------------------------------------------------------
$q2 = $wb.Worksheets.Item("2022-Q2")
$q3 = $wb.Worksheets.Add($q2)
$q3.Name = "2022-Q3"

# ---------------------------------------------------------------------
# 2. Populate "2022-Q3" with the fund-holding table (same shape as the
#    other quarterly sheets: header row + 4 data rows, columns A..H).
# ---------------------------------------------------------------------
$q3.Range("B1").Value = "基金代码"
$q3.Range("C1").Value = "基金名称"
$q3.Range("D1").Value = "基金规模"
$q3.Range("E1").Value = "股票总仓位"
$q3.Range("F1").Value = "仓位占比"
$q3.Range("G1").Value = "持有市值(亿元)"
$q3.Range("H1").Value = "仓位排名"
$q3.Range("B1:H1").Font.Bold = $true
$q3.Range("B1:H1").HorizontalAlignment = -4108
$q3.Range("B1:H1").VerticalAlignment = -4160
$q3.Range("B1:H1").Borders.LineStyle = 1

$q3.Range("A2").Value = 0
$q3.Range("B2").Value = "516150"
$q3.Range("C2").Value = "嘉实中证稀土产业ETF"
$q3.Range("D2").Value = "20.00"
$q3.Range("E2").Value = "99.35"
$q3.Range("F2").Value = "3.81"
$q3.Range("G2").Value = "0.7620"
$q3.Range("H2").Value = 9

$q3.Range("A3").Value = 1
$q3.Range("B3").Value = "516780"
$q3.Range("C3").Value = "华泰柏瑞中证稀土产业ETF"
$q3.Range("D3").Value = "8.05"
$q3.Range("E3").Value = "98.92"
$q3.Range("F3").Value = "3.80"
$q3.Range("G3").Value = "0.3059"
$q3.Range("H3").Value = 9

$q3.Range("A4").Value = 2
$q3.Range("B4").Value = "159715"
$q3.Range("C4").Value = "易方达中证稀土产业ETF"
$q3.Range("D4").Value = "2.52"
$q3.Range("E4").Value = "98.35"
$q3.Range("F4").Value = "3.73"
$q3.Range("G4").Value = "0.0940"
$q3.Range("H4").Value = 9

$q3.Range("A5").Value = 3
$q3.Range("B5").Value = "159713"
$q3.Range("C5").Value = "富国中证稀土产业ETF"
$q3.Range("D5").Value = "2.25"
$q3.Range("E5").Value = "98.40"
$q3.Range("F5").Value = "3.76"
$q3.Range("G5").Value = "0.0846"
$q3.Range("H5").Value = 9

$q3.Range("A2:A5").Font.Bold = $true
$q3.Range("A2:A5").HorizontalAlignment = -4108
$q3.Range("A2:A5").VerticalAlignment = -4160
$q3.Range("A2:A5").Borders.LineStyle = 1

# ---------------------------------------------------------------------
# 3. Insert a new row for "2022-Q3" at the top of the "总计" summary
#    sheet, pushing the older quarters down by one row.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Range("A2:D8").Value = $total.Range("A2:D8").Value
$total.Rows.Item(2).Insert()

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 4
$total.Range("D2").Value = 1.25

$total.Range("A2").Font.Bold = $true
$total.Range("A2").HorizontalAlignment = -4108
$total.Range("A2").VerticalAlignment = -4160
$total.Range("A2").Borders.LineStyle = 1

$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3
$total.Range("A6").Value = 4
$total.Range("A7").Value = 5
$total.Range("A8").Value = 6
